$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need NumberFormat forced to
# Text ("@") first, otherwise Excel (like this COM layer) auto-converts the
# typed text into a floating point Number (losing the original formatted text,
# e.g. "9.20" -> 9.2) exactly as it would for manual keyboard entry.

$ws.Range("D2").Value = "63.702.73"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "2.476.04"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.33"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.94"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").Value = "2.473.02"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.35"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "2.926.16"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "63.567.57"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "2.469.92"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.51"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("E20").Value = "  +6.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.14"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.12"
$ws.Range("E23").Value = "  +18.89%  "
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.08"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "628.87"
$ws.Range("E26").Value = "  +11.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").Value = "  +5.39%  "
$ws.Range("E28").Value = "  +3.54%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.54"
$ws.Range("E29").Value = "  +5.89%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "2.603.39"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.42"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.996"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  +8.27%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.52"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.89"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  +13.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.75"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "151.41"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.61"
$ws.Range("E47").Value = "  +5.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0543"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.608"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0237"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0920"
$ws.Range("E51").Value = "  -0.89%  "
